$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 7
    3  = 1
    4  = 7
    5  = 5
    6  = 4
    7  = 5
    8  = 0
    9  = 5
    10 = 1
    11 = 1
    12 = 1
    13 = 2
    14 = 3
    15 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
